# Apply weekly report update:
#  - refresh "Report Generated On" timestamp
#  - zero out the billed amount summary and bump the line-item count
#  - insert a new line item (Point 17 / GND-MD) into the Monday table,
#    pushing the remaining rows (including the Thursday table) down by one
#  - zero out every "Pricing" (H column) figure in both day tables and
#    their TOTAL rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary updates -------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:48 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 51

# --- Zero out Monday table pricing (rows 16-39) before the insert -------------
for ($r = 16; $r -le 39; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# --- Insert the new "Point 17 / GND-MD" line item as row 40 -------------------
# This pushes the old row 40 ("Point 18") down to row 41, the Monday TOTAL
# down to row 42, and everything in the Thursday table down by one row too.
$ws.Rows("40:40").Insert()

# Excel's row-insert leaves the new blank row 40 with the style copied from
# the row above (the "even" group, 12/13/14) and the shifted-down row 41
# keeps its original "odd" group (9/10/11) style - i.e. the two rows end up
# with their style groups swapped relative to the alternating odd/even
# pattern used throughout the table. Re-sync both rows' formatting from
# unaffected rows above that already carry the correct group so the
# alternating pattern (row40=odd/9-10-11, row41=even/12-13-14) holds.
$ws.Range("A38:H38").Copy()
$ws.Range("A40:H40").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A39:H39").Copy()
$ws.Range("A41:H41").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(40, 1).Value = "Point 17"
$ws.Cells.Item(40, 2).Value = "GND-MD"
$ws.Cells.Item(40, 3).Value = "Inst"
$ws.Cells.Item(40, 4).Value = "GND,Wire Mldg Only"
$ws.Cells.Item(40, 5).Value = "EA"
$ws.Cells.Item(40, 6).Value = 2
$ws.Cells.Item(40, 8).Value = 0

# --- The old "Point 18" row, now row 41, and the Monday TOTAL, now row 42 -----
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(42, 8).Value = 0

# --- Zero out Thursday table pricing (now rows 47-72, shifted down by one) ----
for ($r = 47; $r -le 72; $r++) {
    $ws.Cells.Item($r, 8).Value = 0
}

# --- Nudge the used range back out to column I (touched indirectly by the
#     D:I / G:I merges) so the sheet's reported dimension stays A2:I72 -------
$ws.Range("I72").Font.Bold = $false
